$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

for ($r = 2; $r -le 37; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}
